$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.278.29'
$ws.Range('E2').Value = '  +1.51%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.423.32'
$ws.Range('E3').Value = '  +1.98%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '562.57'
$ws.Range('E5').Value = '  +1.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.10'
$ws.Range('E6').Value = '  +3.17%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +1.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.421.05'
$ws.Range('E9').Value = '  +1.84%  '
$ws.Range('E10').Value = '  +1.51%  '
$ws.Range('E11').Value = '  -2.17%  '
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('E13').Value = '  +0.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.91'
$ws.Range('E14').Value = '  +1.54%  '
$ws.Range('E15').Value = '  +3.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.861.90'
$ws.Range('E16').Value = '  +1.97%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.098.81'
$ws.Range('E17').Value = '  +1.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.423.69'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.31'
$ws.Range('E19').Value = '  +2.71%  '
$ws.Range('E20').Value = '  +1.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.55'
$ws.Range('E21').Value = '  +0.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.74'
$ws.Range('E22').Value = '  +0.64%  '
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('E24').Value = '  +2.18%  '
$ws.Range('E25').Value = '  -3.04%  '
$ws.Range('E26').Value = '  +0.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '582.10'
$ws.Range('E27').Value = '  +11.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.542.95'
$ws.Range('E28').Value = '  +1.96%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0941'
$ws.Range('E30').Value = '  +4.27%  '
$ws.Range('E31').Value = '  +4.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.25'
$ws.Range('E32').Value = '  +0.74%  '
$ws.Range('E33').Value = '  +0.77%  '
$ws.Range('E34').Value = '  +2.47%  '
$ws.Range('E35').Value = '  +1.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.71'
$ws.Range('E36').Value = '  +3.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.78'
$ws.Range('E38').Value = '  +1.86%  '
$ws.Range('E39').Value = '  +1.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '153.00'
$ws.Range('E40').Value = '  +3.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.66'
$ws.Range('E41').Value = '  +0.85%  '
$ws.Range('E42').Value = '  -3.88%  '
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.32'
$ws.Range('E44').Value = '  +7.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '150.19'
$ws.Range('E45').Value = '  +1.74%  '
$ws.Range('E46').Value = '  +1.47%  '
$ws.Range('E47').Value = '  +2.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.31'
$ws.Range('E48').Value = '  +2.78%  '
$ws.Range('E49').Value = '  +2.17%  '
$ws.Range('E50').Value = '  +1.80%  '
$ws.Range('E51').Value = '  +1.84%  '
